$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update financial figures for rows 2-6 (2014/12 .. 2018/12 IFRS-consolidated)
# Row 2
$ws.Range("D2").Value = 5851
$ws.Range("E2").Value = 692
$ws.Range("F2").Value = 692
$ws.Range("G2").Value = 641
$ws.Range("H2").Value = 483
$ws.Range("I2").Value = 391
$ws.Range("J2").Value = 92
$ws.Range("K2").Value = 6172
$ws.Range("L2").Value = 2830
$ws.Range("M2").Value = 3343
$ws.Range("N2").Value = 2560
$ws.Range("O2").Value = 783
$ws.Range("P2").Value = 88
$ws.Range("Q2").Value = 356
$ws.Range("R2").Value = -217
$ws.Range("S2").Value = 18
$ws.Range("T2").Value = 186
$ws.Range("U2").Value = 169
$ws.Range("V2").Value = 1567
$ws.Range("W2").Value = 11.83
$ws.Range("X2").Value = 8.26
$ws.Range("Y2").Value = 16.45
$ws.Range("Z2").Value = 8.16
$ws.Range("AA2").Value = 84.65
$ws.Range("AB2").Value = 4011.02
$ws.Range("AC2").Value = 1961
$ws.Range("AD2").Value = 5.88
$ws.Range("AE2").Value = 13284
$ws.Range("AF2").Value = 0.87
$ws.Range("AG2").Value = 53
$ws.Range("AH2").Value = 0.46
$ws.Range("AI2").Value = 2.59
$ws.Range("AJ2").Value = 19930000

# Row 3
$ws.Range("D3").Value = 6031
$ws.Range("E3").Value = 737
$ws.Range("F3").Value = 737
$ws.Range("G3").Value = 736
$ws.Range("H3").Value = 526
$ws.Range("I3").Value = 412
$ws.Range("J3").Value = 113
$ws.Range("K3").Value = 6391
$ws.Range("L3").Value = 2540
$ws.Range("M3").Value = 3851
$ws.Range("N3").Value = 2958
$ws.Range("O3").Value = 894
$ws.Range("P3").Value = 90
$ws.Range("Q3").Value = 635
$ws.Range("R3").Value = -365
$ws.Range("S3").Value = -135
$ws.Range("T3").Value = 413
$ws.Range("U3").Value = 222
$ws.Range("V3").Value = 1433
$ws.Range("W3").Value = 12.23
$ws.Range("X3").Value = 8.72
$ws.Range("Y3").Value = 14.95
$ws.Range("Z3").Value = 8.37
$ws.Range("AA3").Value = 65.94
$ws.Range("AB3").Value = 4343.09
$ws.Range("AC3").Value = 2069
$ws.Range("AD3").Value = 5.81
$ws.Range("AE3").Value = 15346
$ws.Range("AF3").Value = 0.78
$ws.Range("AG3").Value = 54
$ws.Range("AH3").Value = 0.45
$ws.Range("AI3").Value = 2.52
$ws.Range("AJ3").Value = 19930000

# Row 4
$ws.Range("D4").Value = 5513
$ws.Range("E4").Value = 675
$ws.Range("F4").Value = 675
$ws.Range("G4").Value = 666
$ws.Range("H4").Value = 489
$ws.Range("I4").Value = 394
$ws.Range("J4").Value = 95
$ws.Range("K4").Value = 7012
$ws.Range("L4").Value = 2708
$ws.Range("M4").Value = 4304
$ws.Range("N4").Value = 3317
$ws.Range("O4").Value = 987
$ws.Range("P4").Value = 93
$ws.Range("Q4").Value = 623
$ws.Range("R4").Value = -664
$ws.Range("S4").Value = 107
$ws.Range("T4").Value = 633
$ws.Range("U4").Value = -9
$ws.Range("V4").Value = 1542
$ws.Range("W4").Value = 12.24
$ws.Range("X4").Value = 8.87
$ws.Range("Y4").Value = 12.55
$ws.Range("Z4").Value = 7.3
$ws.Range("AA4").Value = 62.92
$ws.Range("AB4").Value = 4639.13
$ws.Range("AC4").Value = 1976
$ws.Range("AD4").Value = 5.96
$ws.Range("AE4").Value = 17210
$ws.Range("AF4").Value = 0.68
$ws.Range("AG4").Value = 74
$ws.Range("AH4").Value = 0.63
$ws.Range("AI4").Value = 3.63
$ws.Range("AJ4").Value = 19930000

# Row 5
$ws.Range("D5").Value = 5439
$ws.Range("E5").Value = 372
$ws.Range("F5").Value = 372
$ws.Range("G5").Value = 354
$ws.Range("H5").Value = 191
$ws.Range("I5").Value = 180
$ws.Range("J5").Value = 12
$ws.Range("K5").Value = 6774
$ws.Range("L5").Value = 2404
$ws.Range("M5").Value = 4371
$ws.Range("N5").Value = 3383
$ws.Range("O5").Value = 987
$ws.Range("P5").Value = 95
$ws.Range("Q5").Value = 221
$ws.Range("R5").Value = -412
$ws.Range("S5").Value = -73
$ws.Range("T5").Value = 502
$ws.Range("U5").Value = -281
$ws.Range("V5").Value = 1422
$ws.Range("W5").Value = 6.83
$ws.Range("X5").Value = 3.52
$ws.Range("Y5").Value = 5.36
$ws.Range("Z5").Value = 2.77
$ws.Range("AA5").Value = 54.99
$ws.Range("AB5").Value = 4682.51
$ws.Range("AC5").Value = 902
$ws.Range("AD5").Value = 8.52
$ws.Range("AE5").Value = 17555
$ws.Range("AF5").Value = 0.44
$ws.Range("AG5").Value = 95
$ws.Range("AH5").Value = 1.24
$ws.Range("AI5").Value = 10.21
$ws.Range("AJ5").Value = 19930000

# Row 6
$ws.Range("D6").Value = 5454
$ws.Range("E6").Value = 270
$ws.Range("F6").Value = 270
$ws.Range("G6").Value = 229
$ws.Range("H6").Value = 110
$ws.Range("I6").Value = 59
$ws.Range("K6").Value = 7142
$ws.Range("L6").Value = 2680
$ws.Range("M6").Value = 4463
$ws.Range("N6").Value = 3425
$ws.Range("P6").Value = 97
$ws.Range("Q6").Value = 278
$ws.Range("R6").Value = -450
$ws.Range("S6").Value = 194
$ws.Range("T6").Value = 487
$ws.Range("U6").Value = -209
$ws.Range("V6").Value = 1620
$ws.Range("W6").Value = 4.96
$ws.Range("X6").Value = 2.01
$ws.Range("Y6").Value = 1.73
$ws.Range("Z6").Value = 1.58
$ws.Range("AA6").Value = 60.04
$ws.Range("AB6").Value = 4625.05
$ws.Range("AC6").Value = 296
$ws.Range("AD6").Value = 22.16
$ws.Range("AE6").Value = 17770
$ws.Range("AF6").Value = 0.37
$ws.Range("AG6").Value = 108
$ws.Range("AH6").Value = 1.64
$ws.Range("AI6").Value = 35.09
$ws.Range("AJ6").Value = 19930000

# Remove stale estimate data (2019/12(E) .. 2021/12(E)) - keep only rank/name columns
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()

Write-Host "Applied IFRS list corrections"